$wb = $excel.ActiveWorkbook

# --- Rename the two survey sheets for parallel structure with the notes sheets ---
$wsSurvey1880 = $wb.Worksheets.Item("1881")
$wsSurvey1880.Name = "1880Survey"

$wsSurvey1940 = $wb.Worksheets.Item("1940")
$wsSurvey1940.Name = "1940Survey"

# --- Add a new metadata sheet at the very end explaining the naming convention ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$metaSheet = $wb.Worksheets.Add($null, $lastSheet)
$metaSheet.Name = "1880Metadata"
$metaCell = $metaSheet.Range("A1")
$metaCell.Value = "Actually surveyed in 1881; tab reads 1880 for consistency in the data processing R script"
$metaCell.Font.Color = 0

# --- Scroll the 1940Survey sheet view down (best effort) ---
[void]$wsSurvey1940.Activate()
$excel.ActiveWindow.ScrollRow = 54
$excel.ActiveWindow.ScrollColumn = 1
[void]$wsSurvey1940.Range("G2:G79").Select()

# --- Make the 1880Survey sheet the active/selected tab ---
[void]$wsSurvey1880.Activate()
